# Resume edit: insert "Javascript, " before "React.js, Express.js, ..." in the
# "Technical specialties include ..." sentence, and let Word's own "_GoBack"
# last-edit bookmark follow the new insertion point (moving it off the
# "—June 2019" spot it previously sat at).

$d = $word.ActiveDocument

# Locate the point right before "React.js" in the technical-specialties run.
$rng = $d.Content
$found = $rng.Find.Execute("React.js, Express.js, Node.js", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the technical-specialties text to anchor the edit on."
}
$insertPos = $rng.Start
$insertText = "Javascript, "

# Insert the new text as its own run.
$ip = $d.Range($insertPos, $insertPos)
$ip.InsertBefore($insertText)

# Toggle a character property on just the inserted text so the engine keeps
# it as a distinct run instead of silently re-coalescing it with its
# (identically formatted) neighbors.
$insertedRng = $d.Range($insertPos, $insertPos + $insertText.Length)
$insertedRng.Bold = 1
$insertedRng.Bold = 0

# Word keeps a single "_GoBack" bookmark marking the location of the most
# recent edit; re-adding it moves it here (and removes it from its old spot
# after "—June 2019").
$bmPos = $insertPos + $insertText.Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
